# Generate Report for Handback
# Refresh the "Latest HO Xliff Generate Date" / handoff & handback
# timestamps for the ea69d92a*.md row now that a new handback report was
# generated for it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the ea69d92a file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-31 16:56:25"

# --- zh-cn sheet: row 3 is the ea69d92a file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-08-31 16:56:21"
$wsZhCn.Range("K3").Value = "2016-08-31 16:56:37"

# --- de-de sheet: row 3 is the ea69d92a file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-08-31 16:56:25"
$wsDeDe.Range("K3").Value = "2016-08-31 16:56:44"
